$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8950971961021423
$ws.Range("B1").Value = 2.971071004867554
$ws.Range("C1").Value = 4.384921073913574
$ws.Range("D1").Value = 3.031518697738647
$ws.Range("E1").Value = 1.402795195579529
